$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the team-lookup table (columns H:J)
$ws.Range("H1").Value = "Name"
$ws.Range("I1").Value = "ID"
$ws.Range("J1").Value = "Acronym"

# Team acronyms, one per lookup row (H2:H21 already hold the team names,
# I2:I21 their numeric IDs) - fill in column J with the short acronym.
$ws.Range("J2").Value  = "YB"
$ws.Range("J3").Value  = "BAS"
$ws.Range("J4").Value  = "SFC"
$ws.Range("J5").Value  = "LUG"
$ws.Range("J6").Value  = "LUZ"
$ws.Range("J7").Value  = "LS"
$ws.Range("J8").Value  = "SG"
$ws.Range("J9").Value  = "FCZ"
$ws.Range("J10").Value = "SIO"
$ws.Range("J11").Value = "GC"
$ws.Range("J12").Value = "VAD"
$ws.Range("J13").Value = "THU"
$ws.Range("J14").Value = "SLO"
$ws.Range("J15").Value = "FCS"
$ws.Range("J16").Value = "AAR"
$ws.Range("J17").Value = "WIN"
$ws.Range("J18").Value = "WILL"
$ws.Range("J19").Value = "SCK"
$ws.Range("J20").Value = "XAM"
$ws.Range("J21").Value = "YS"

# Leave the view scrolled/selected near the new table, like the author did
# while reviewing the freshly entered acronyms.
$null = $ws.Range("K11").Select()
